$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''D...
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.5, early_stopping_rounds=3,
                                learning_rate=0.05, max_bin=75, max_depth=7,
                                min_child_samples=7, min_data_in_leaf=25,
                                num_iterations=400, num_leaves=2,
                                random_state=42, subsample=0.5))])'
$ws.Range("C2").Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__num_iterations'': 400, ''model__min_data_in_leaf'': 25, ''model__min_child_samples'': 7, ''model__max_depth'': 7, ''model__max_bin'': 75, ''model__learning_rate'': 0.05, ''model__early_stopping_rounds'': 3, ''model__colsample_bytree'': 0.5, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D2").Value = 0.5887404140631359
$ws.Range("G2").Value = 0.8279569892473119
$ws.Range("H2").Value = 0.539241622574956
$ws.Range("I2").Value = '[1 0 1 0 0 0 1 1 1 1 1 1 1 0 0 0 0 0 0 0 1 1 0 0]'
$ws.Range("J2").Value = '[0 0 0 0 0 1 0 1 1 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1]'

$ws.Range("B3").Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''D...
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.5, early_stopping_rounds=7,
                                learning_rate=0.05, max_bin=25, max_depth=7,
                                min_child_samples=5, min_data_in_leaf=30,
                                num_iterations=300, num_leaves=2,
                                random_state=42, subsample=0.5))])'
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__num_iterations'': 300, ''model__min_data_in_leaf'': 30, ''model__min_child_samples'': 5, ''model__max_depth'': 7, ''model__max_bin'': 25, ''model__learning_rate'': 0.05, ''model__early_stopping_rounds'': 7, ''model__colsample_bytree'': 0.5, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D3").Value = 0.5786331989976492
$ws.Range("G3").Value = 0.7420548476619618
$ws.Range("H3").Value = 0.4999999999999998
$ws.Range("I3").Value = '[0 1 1 0 1 0 0 0 1 1 1 0 0 0 1 0 1 0 1 1 0 0 1 0]'
$ws.Range("J3").Value = '[0 1 1 1 0 0 0 0 1 1 1 1 1 1 0 1 0 0 1 0 0 1 0 1]'

$ws.Range("B4").Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''D...
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.4, early_stopping_rounds=1,
                                max_bin=50, max_depth=7, min_child_samples=7,
                                min_data_in_leaf=25, num_iterations=300,
                                num_leaves=2, random_state=42,
                                subsample=0.5))])'
$ws.Range("C4").Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__num_iterations'': 300, ''model__min_data_in_leaf'': 25, ''model__min_child_samples'': 7, ''model__max_depth'': 7, ''model__max_bin'': 50, ''model__learning_rate'': 0.1, ''model__early_stopping_rounds'': 1, ''model__colsample_bytree'': 0.4, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D4").Value = 0.6458327413900169
$ws.Range("G4").Value = 0.8495319944694777
$ws.Range("H4").Value = 0.3584656084656084
$ws.Range("I4").Value = '[0 1 1 0 0 1 0 0 0 1 0 1 0 1 0 1 0 0 1 0 0 1 1 1]'
$ws.Range("J4").Value = '[1 0 1 1 1 1 1 1 1 0 1 0 1 1 0 1 0 1 1 0 1 0 0 1]'

$ws.Range("B5").Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''D...
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.5, early_stopping_rounds=7,
                                learning_rate=0.05, max_bin=50, max_depth=7,
                                min_child_samples=7, min_data_in_leaf=30,
                                num_iterations=300, num_leaves=2,
                                random_state=42, subsample=0.5))])'
$ws.Range("C5").Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__num_iterations'': 300, ''model__min_data_in_leaf'': 30, ''model__min_child_samples'': 7, ''model__max_depth'': 7, ''model__max_bin'': 50, ''model__learning_rate'': 0.05, ''model__early_stopping_rounds'': 7, ''model__colsample_bytree'': 0.5, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D5").Value = 0.5916797340945947
$ws.Range("G5").Value = 0.7850457063849681
$ws.Range("H5").Value = 0.4928571428571429
$ws.Range("I5").Value = '[0 0 0 1 1 1 1 1 0 0 1 0 1 0 0 0 1 0 1 1 0 0 1 0]'
$ws.Range("J5").Value = '[1 0 1 0 0 1 1 1 0 0 0 1 1 0 0 1 0 1 0 0 0 0 0 0]'

$ws.Range("B6").Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''D...
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.5, early_stopping_rounds=7,
                                max_bin=25, max_depth=7, min_child_samples=5,
                                min_data_in_leaf=25, num_iterations=400,
                                num_leaves=2, random_state=42,
                                subsample=0.5))])'
$ws.Range("C6").Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__num_iterations'': 400, ''model__min_data_in_leaf'': 25, ''model__min_child_samples'': 5, ''model__max_depth'': 7, ''model__max_bin'': 25, ''model__learning_rate'': 0.1, ''model__early_stopping_rounds'': 7, ''model__colsample_bytree'': 0.5, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D6").Value = 0.5600832629702598
$ws.Range("G6").Value = 0.8922980745186295
$ws.Range("H6").Value = 0.5392753623188405
$ws.Range("I6").Value = '[1 0 1 0 1 1 0 0 0 1 0 1 0 1 1 1 0 1 0 1 0 0 0 0]'
$ws.Range("J6").Value = '[0 0 0 1 1 0 0 0 1 0 1 1 0 1 1 1 0 1 0 1 1 1 1 1]'
